# CharacterDataSheet.xlsx - renumber character IDs in column A to start
# at 0 instead of 1 (rows 2-30), and move the sheet selection down to A31
# (the first empty row below the table), matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CharacterDatas")

for ($r = 2; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    $cell.Value = $current - 1
}

$ws.Range("A31").Select()
